$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header "Save" in H1, matching the style/formatting of the existing header row (G1)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add data values for the new "Save" column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
